$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 276.7
$ws.Range("I11").Value = 276.7
$ws.Range("K11").Value = 276.7
$ws.Range("M11").Value = -136.7
$ws.Range("H125").Value = 1808.8
$ws.Range("I125").Value = 432
$ws.Range("J125").Value = 2726.6667
$ws.Range("K125").Value = 3888
$ws.Range("L125").Value = 24540.0003
$ws.Range("M125").Value = -1428
$ws.Range("N125").Value = -29460.0003
$ws.Range("H129").Value = 842.4737
$ws.Range("I129").Value = 679.4
$ws.Range("J129").Value = 858.1539
$ws.Range("K129").Value = 2038.2
$ws.Range("L129").Value = 2574.4617
$ws.Range("M129").Value = 2961.8
$ws.Range("N129").Value = -12574.4617
$ws.Range("H137").Value = 45536.434
$ws.Range("I137").Value = 2081.0625
$ws.Range("J137").Value = 144863
$ws.Range("K137").Value = 6243.1875
$ws.Range("L137").Value = 434589
$ws.Range("M137").Value = -3693.1875
$ws.Range("N137").Value = -439689

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 20023.418
$ws.Range("I32").Value = 22879.064
$ws.Range("J32").Value = 3246.5
$ws.Range("K32").Value = 22879.064
$ws.Range("L32").Value = 3246.5
$ws.Range("M32").Value = -22592.064
$ws.Range("N32").Value = -3820.5
$ws.Range("H45").Value = 3917.08
$ws.Range("I45").Value = 4559
$ws.Range("K45").Value = 4559
$ws.Range("M45").Value = -4182
$ws.Range("H97").Value = 1939.7693
$ws.Range("I97").Value = 1671.3334
$ws.Range("J97").Value = 2543.75
$ws.Range("K97").Value = 1671.3334
$ws.Range("L97").Value = 2543.75
$ws.Range("M97").Value = -1175.3334
$ws.Range("N97").Value = -3535.75
$ws.Range("H122").Value = 1835.6364
$ws.Range("I122").Value = 1919.2
$ws.Range("K122").Value = 5757.6
$ws.Range("M122").Value = -3307.6
$ws.Range("H127").Value = 49390
$ws.Range("J127").Value = 49390
$ws.Range("L127").Value = 49390
$ws.Range("N127").Value = -59310

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1962.5
$ws.Range("I86").Value = 1794.3636
$ws.Range("K86").Value = 1794.3636
$ws.Range("M86").Value = -671.3635999999999
$ws.Range("H89").Value = 1962.5
$ws.Range("I89").Value = 1794.3636
$ws.Range("K89").Value = 8971.817999999999
$ws.Range("M89").Value = -3355.817999999999
$ws.Range("H105").Value = 2825
$ws.Range("I105").Value = 3266.6667
$ws.Range("J105").Value = 2560
$ws.Range("K105").Value = 3266.6667
$ws.Range("L105").Value = 2560
$ws.Range("M105").Value = -1519.6667
$ws.Range("N105").Value = -6054
$ws.Range("H119").Value = 19879.5
$ws.Range("J119").Value = 19879.5
$ws.Range("L119").Value = 19879.5
$ws.Range("N119").Value = -29555.5
$ws.Range("H132").Value = 49995
$ws.Range("J132").Value = 49995
$ws.Range("L132").Value = 49995
$ws.Range("N132").Value = -60115
$ws.Range("H137").Value = 50740
$ws.Range("J137").Value = 50740
$ws.Range("L137").Value = 50740
$ws.Range("N137").Value = -60940

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1339
$ws.Range("I16").Value = 1298.75
$ws.Range("J16").Value = 1500
$ws.Range("K16").Value = 1298.75
$ws.Range("L16").Value = 1500
$ws.Range("M16").Value = -1011.75
$ws.Range("N16").Value = -2074
$ws.Range("H20").Value = 46999.6
$ws.Range("J20").Value = 46999.6
$ws.Range("L20").Value = 46999.6
$ws.Range("N20").Value = -47471.6
$ws.Range("H30").Value = 46999.6
$ws.Range("J30").Value = 46999.6
$ws.Range("L30").Value = 46999.6
$ws.Range("N30").Value = -47181.6
$ws.Range("H31").Value = 12527.611
$ws.Range("I31").Value = 15650.038
$ws.Range("J31").Value = 4409.3
$ws.Range("K31").Value = 15650.038
$ws.Range("L31").Value = 4409.3
$ws.Range("M31").Value = -15355.038
$ws.Range("N31").Value = -4999.3
$ws.Range("H34").Value = 12527.611
$ws.Range("I34").Value = 15650.038
$ws.Range("J34").Value = 4409.3
$ws.Range("K34").Value = 15650.038
$ws.Range("L34").Value = 4409.3
$ws.Range("M34").Value = -15448.038
$ws.Range("N34").Value = -4813.3
$ws.Range("H58").Value = 21897.541
$ws.Range("I58").Value = 1348.6
$ws.Range("J58").Value = 56145.777
$ws.Range("K58").Value = 1348.6
$ws.Range("L58").Value = 56145.777
$ws.Range("M58").Value = -1145.6
$ws.Range("N58").Value = -56551.777
$ws.Range("H113").Value = 1339
$ws.Range("I113").Value = 1298.75
$ws.Range("J113").Value = 1500
$ws.Range("K113").Value = 1298.75
$ws.Range("L113").Value = 1500
$ws.Range("M113").Value = 871.25
$ws.Range("N113").Value = -5840
$ws.Range("H128").Value = 46999.6
$ws.Range("J128").Value = 46999.6
$ws.Range("L128").Value = 46999.6
$ws.Range("N128").Value = -56959.6
$ws.Range("H136").Value = 21897.541
$ws.Range("I136").Value = 1348.6
$ws.Range("J136").Value = 56145.777
$ws.Range("K136").Value = 4045.8
$ws.Range("L136").Value = 168437.331
$ws.Range("M136").Value = -1495.8
$ws.Range("N136").Value = -173537.331

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1213.7858
$ws.Range("I5").Value = 1008.0909
$ws.Range("K5").Value = 3024.2727
$ws.Range("M5").Value = -2912.2727
$ws.Range("H7").Value = 0
$ws.Range("I7").Value = 0
$ws.Range("K7").Value = 0
$ws.Range("M7").ClearContents()
$ws.Range("H9").Value = 1001
$ws.Range("J9").Value = 1001
$ws.Range("L9").Value = 3003
$ws.Range("N9").Value = -3451
$ws.Range("H33").Value = 228.66667
$ws.Range("I33").Value = 193
$ws.Range("K33").Value = 1158
$ws.Range("M33").Value = -875
$ws.Range("H54").Value = 4004.1667
$ws.Range("J54").Value = 4004.1667
$ws.Range("L54").Value = 12012.5001
$ws.Range("N54").Value = -13130.5001
$ws.Range("H117").Value = 37038880
$ws.Range("I117").Value = 829.6
$ws.Range("J117").Value = 83336450
$ws.Range("K117").Value = 2488.8
$ws.Range("L117").Value = 250009350
$ws.Range("M117").Value = 953.1999999999998
$ws.Range("N117").Value = -250016234
$ws.Range("H131").Value = 111915.05
$ws.Range("J131").Value = 115739.36
$ws.Range("L131").Value = 347218.08
$ws.Range("N131").Value = -357298.08
$ws.Range("H132").Value = 794.25
$ws.Range("I132").Value = 794.25
$ws.Range("K132").Value = 7148.25
$ws.Range("M132").Value = -4618.25
$ws.Range("H135").Value = 1213.7858
$ws.Range("I135").Value = 1008.0909
$ws.Range("K135").Value = 9072.8181
$ws.Range("M135").Value = -6537.8181

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2207.3928
$ws.Range("I102").Value = 2315.9167
$ws.Range("J102").Value = 1556.25
$ws.Range("K102").Value = 2315.9167
$ws.Range("L102").Value = 1556.25
$ws.Range("M102").Value = -693.9167000000002
$ws.Range("N102").Value = -4800.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 248.33333
$ws.Range("I55").Value = 137.625
$ws.Range("J55").Value = 469.75
$ws.Range("K55").Value = 137.625
$ws.Range("L55").Value = 469.75
$ws.Range("M55").Value = 35.375
$ws.Range("N55").Value = -815.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 2540
$ws.Range("H132").Value = 3045.923
$ws.Range("I132").Value = 2690.818
$ws.Range("J132").Value = 4999
$ws.Range("K132").Value = 8072.454000000001
$ws.Range("L132").Value = 14997
$ws.Range("M132").Value = -5542.454000000001
$ws.Range("N132").Value = -20057
$ws.Range("H136").Value = 1110.7106
$ws.Range("I136").Value = 842.4231
$ws.Range("J136").Value = 1692
$ws.Range("K136").Value = 2527.2693
$ws.Range("L136").Value = 5076
$ws.Range("M136").Value = 22.73070000000007
$ws.Range("N136").Value = -10176

Write-Output "Applied market data refresh to 8 sheets"